$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new data to existing row 5 (X5, Y5)
$ws.Range("X5").Value = 0.6499990000000011
$ws.Range("Y5").Value = "Up"

# Add new row 6 with a fresh scan result
$ws.Range("A5").Copy($ws.Range("A6"))
$ws.Range("A6").Value = 42647.884270833332
$ws.Range("B6").Value = 11
$ws.Range("C6").Value = "Buy"
$ws.Range("D6").Value = 24
$ws.Range("E6").Value = 14585
$ws.Range("F6").Value = 2144
$ws.Range("G6").Value = 61
$ws.Range("H6").Value = 35
$ws.Range("I6").Value = 82
$ws.Range("J6").Value = 17
$ws.Range("K6").Value = 24061
$ws.Range("L6").Value = 343
$ws.Range("M6").Value = 196
$ws.Range("N6").Value = 104
$ws.Range("O6").Value = 22
$ws.Range("P6").Value = "Noun"
$ws.Range("Q6").Value = 46.357611069683557
$ws.Range("R6").Value = 0
$ws.Range("S5").Copy($ws.Range("S6"))
$ws.Range("S6").Value = 0.0591
$ws.Range("T5").Copy($ws.Range("T6"))
$ws.Range("T6").Value = -0.0421
$ws.Range("U6").Value = 2.25
$ws.Range("V6").Value = "N/A"
$ws.Range("W6").Value = 0
